# Daily attendance processing - 2026-01-22 19:18:22
# Re-order the comma-separated "Recorded By" author lists in column G.
# The author set per cell is unchanged; only the ordering of the
# names/emails within the string changes (matches the source diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact before -> after string replacements observed for the
# "Recorded By" column (G). Applied wherever the current cell text
# matches one of the "before" values.
$map = @{
    'system, backup@backdoor.com, System' = 'system, System, backup@backdoor.com'
    'backup@backdoor.com, System'         = 'System, backup@backdoor.com'
    'admin@admin.com, System'             = 'System, admin@admin.com'
    'System, dnasr281@gmail.com'          = 'dnasr281@gmail.com, System'
    'admin@admin.com, dnasr281@gmail.com' = 'dnasr281@gmail.com, admin@admin.com'
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value2 = $map[$val]
    }
}
